$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.564.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.52%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.389.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.71%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.32%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'406.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.54%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'125.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.79%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.611"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.85%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.709"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -4.16%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.131"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -9.80%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'41.70"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.32%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -0.33%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.921.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.81%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'8.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.94%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'20.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.63%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0000200"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -8.02%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.384.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.08%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'12.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.49%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -0.52%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'61.569.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.48%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'478.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +20.03%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'88.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.75%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'3.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.15%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'12.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.35%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'3.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.94%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'32.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.41%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +4.99%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.12%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +3.11%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'2.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -4.86%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'11.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.78%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.166"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.09%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -6.32%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'40.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -5.60%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.70%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'54.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.10%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.0477"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -3.89%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.14%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'TheGraph"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'0.325"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +4.81%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'Monero"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'145.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.24%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -2.06%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.132"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.13%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.05%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +3.43%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +4.43%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'4.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.09%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +19.58%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'16.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.87%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Cronos"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.142"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +7.92%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'EnergySwap"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'21.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.87%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'111.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +14.47%  "
$ws.Range("E51").Style = "Normal"
